# Refresh market-price-derived columns (currentAveragePrice.. LeveProfitHQ) for the
# rows whose Universalis quotes moved since the last scheduled run.
# Values below are taken verbatim from the upstream data refresh.
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 78690.55499999999
$ws.Range("I28").Value = 100629.86
$ws.Range("J28").Value = 1903
$ws.Range("K28").Value = 100629.86
$ws.Range("L28").Value = 1903
$ws.Range("M28").Value = -100144.86
$ws.Range("N28").Value = -2873
# Row 86
$ws.Range("H86").Value = 2025.12
$ws.Range("I86").Value = 2351.8667
$ws.Range("K86").Value = 2351.8667
$ws.Range("M86").Value = -1228.8667
# Row 89
$ws.Range("H89").Value = 2025.12
$ws.Range("I89").Value = 2351.8667
$ws.Range("K89").Value = 11759.3335
$ws.Range("M89").Value = -6143.333500000001
# Row 98
$ws.Range("H98").Value = 2700
$ws.Range("I98").Value = 2700
$ws.Range("K98").Value = 2700
$ws.Range("M98").Value = -1202
# Row 113
$ws.Range("H113").Value = 2483
$ws.Range("I113").Value = 1978.75
$ws.Range("K113").Value = 1978.75
$ws.Range("M113").Value = 1275.25
# Row 116
$ws.Range("H116").Value = 5459.8647
$ws.Range("I116").Value = 5215.5356
$ws.Range("J116").Value = 6220
$ws.Range("K116").Value = 5215.5356
$ws.Range("L116").Value = 6220
$ws.Range("M116").Value = -1773.5356
$ws.Range("N116").Value = -13104
# Row 122
$ws.Range("H122").Value = 2700
$ws.Range("I122").Value = 2700
$ws.Range("K122").Value = 8100
$ws.Range("M122").Value = -5650
# Row 132
$ws.Range("H132").Value = 2747.3965
$ws.Range("I132").Value = 2344.425
$ws.Range("J132").Value = 3642.889
$ws.Range("K132").Value = 7033.275000000001
$ws.Range("L132").Value = 10928.667
$ws.Range("M132").Value = -4503.275000000001
$ws.Range("N132").Value = -15988.667
# Row 137
$ws.Range("H137").Value = 3839.2888
$ws.Range("I137").Value = 1136.375
$ws.Range("J137").Value = 5330.552
$ws.Range("K137").Value = 3409.125
$ws.Range("L137").Value = 15991.656
$ws.Range("M137").Value = -859.125
$ws.Range("N137").Value = -21091.656

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 19500
$ws.Range("J24").Value = 19500
$ws.Range("L24").Value = 19500
$ws.Range("N24").Value = -20248
# Row 61
$ws.Range("H61").Value = 33336138
$ws.Range("I61").Value = 52633796
$ws.Range("J61").Value = 3818.3635
$ws.Range("K61").Value = 52633796
$ws.Range("L61").Value = 3818.3635
$ws.Range("M61").Value = -52633584
$ws.Range("N61").Value = -4242.363499999999
# Row 74
$ws.Range("H74").Value = 12671.333
$ws.Range("J74").Value = 16063.143
$ws.Range("L74").Value = 16063.143
$ws.Range("N74").Value = -17811.143
# Row 77
$ws.Range("H77").Value = 12671.333
$ws.Range("J77").Value = 16063.143
$ws.Range("L77").Value = 80315.715
$ws.Range("N77").Value = -89051.715
# Row 98
$ws.Range("H98").Value = 14996.5
$ws.Range("J98").Value = 14996.5
$ws.Range("L98").Value = 14996.5
$ws.Range("N98").Value = -20986.5
# Row 100
$ws.Range("H100").Value = 19500
$ws.Range("J100").Value = 19500
$ws.Range("L100").Value = 19500
$ws.Range("N100").Value = -21664
# Row 132
$ws.Range("H132").Value = 11051.708
$ws.Range("I132").Value = 10125.154
$ws.Range("J132").Value = 12146.728
$ws.Range("K132").Value = 30375.462
$ws.Range("L132").Value = 36440.18399999999
$ws.Range("M132").Value = -27845.462
$ws.Range("N132").Value = -41500.18399999999
# Row 136
$ws.Range("H136").Value = 33336138
$ws.Range("I136").Value = 52633796
$ws.Range("J136").Value = 3818.3635
$ws.Range("K136").Value = 157901388
$ws.Range("L136").Value = 11455.0905
$ws.Range("M136").Value = -157898838
$ws.Range("N136").Value = -16555.0905

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
# Row 100
$ws.Range("H100").Value = 22321.5
$ws.Range("J100").Value = 22321.5
$ws.Range("L100").Value = 22321.5
$ws.Range("N100").Value = -24485.5

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2563.5264
$ws.Range("I16").Value = 1357.2307
$ws.Range("J16").Value = 5177.1665
$ws.Range("K16").Value = 1357.2307
$ws.Range("L16").Value = 5177.1665
$ws.Range("M16").Value = -1070.2307
$ws.Range("N16").Value = -5751.1665
# Row 28
$ws.Range("H28").Value = 48000
$ws.Range("J28").Value = 48000
$ws.Range("L28").Value = 48000
$ws.Range("N28").Value = -48490
# Row 59
$ws.Range("H59").Value = 26750
$ws.Range("J59").Value = 26750
$ws.Range("L59").Value = 26750
$ws.Range("N59").Value = -29040
# Row 106
$ws.Range("H106").Value = 23166.666
$ws.Range("J106").Value = 23166.666
$ws.Range("L106").Value = 23166.666
$ws.Range("N106").Value = -25690.666
# Row 113
$ws.Range("H113").Value = 2563.5264
$ws.Range("I113").Value = 1357.2307
$ws.Range("J113").Value = 5177.1665
$ws.Range("K113").Value = 1357.2307
$ws.Range("L113").Value = 5177.1665
$ws.Range("M113").Value = 812.7692999999999
$ws.Range("N113").Value = -9517.166499999999

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 3583.75
$ws.Range("I80").Value = 2500.3333
$ws.Range("J80").Value = 4667.1665
$ws.Range("K80").Value = 7500.999899999999
$ws.Range("L80").Value = 14001.4995
$ws.Range("M80").Value = -6564.999899999999
$ws.Range("N80").Value = -15873.4995
# Row 83
$ws.Range("H83").Value = 3583.75
$ws.Range("I83").Value = 2500.3333
$ws.Range("J83").Value = 4667.1665
$ws.Range("K83").Value = 22502.9997
$ws.Range("L83").Value = 42004.4985
$ws.Range("M83").Value = -17822.9997
$ws.Range("N83").Value = -51364.4985
# Row 98
$ws.Range("H98").Value = 401.03226
$ws.Range("J98").Value = 479.73334
$ws.Range("L98").Value = 1439.20002
$ws.Range("N98").Value = -4435.20002
# Row 103
$ws.Range("H103").Value = 841.6667
$ws.Range("I103").Value = 116.666664
$ws.Range("K103").Value = 349.999992
$ws.Range("M103").Value = 529.000008
# Row 122
$ws.Range("H122").Value = 1382.8572
$ws.Range("I122").Value = 596
$ws.Range("K122").Value = 5364
$ws.Range("M122").Value = -2914

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
# Row 98
$ws.Range("H98").Value = 9708.6
$ws.Range("J98").Value = 9708.6
$ws.Range("L98").Value = 9708.6
$ws.Range("N98").Value = -15698.6
# Row 113
$ws.Range("H113").Value = 4041.087
$ws.Range("J113").Value = 5085.5
$ws.Range("L113").Value = 5085.5
$ws.Range("N113").Value = -9425.5
# Row 122
$ws.Range("H122").Value = 1370.7778
$ws.Range("I122").Value = 995.6667
$ws.Range("J122").Value = 1558.3334
$ws.Range("K122").Value = 2987.0001
$ws.Range("L122").Value = 4675.0002
$ws.Range("M122").Value = -537.0001000000002
$ws.Range("N122").Value = -9575.0002

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
# Row 99
$ws.Range("H99").Value = 27428.334
$ws.Range("I99").Value = 22500
$ws.Range("J99").Value = 37285
$ws.Range("K99").Value = 22500
$ws.Range("L99").Value = 37285
$ws.Range("M99").Value = -19505
$ws.Range("N99").Value = -43275

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
# Row 98
$ws.Range("H98").Value = 31599.8
$ws.Range("J98").Value = 31599.8
$ws.Range("L98").Value = 31599.8
$ws.Range("N98").Value = -37589.8
